# Weekly fruit/vegetable price update:
# Insert a new record row at row 362 (shifting the existing rows 362..384
# down to 363..385) and populate it with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 362:384 down by inserting a new blank row at 362.
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row 362 with the new observation.
$ws.Cells.Item(362, 1).Value = 7
$ws.Cells.Item(362, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(362, 3).Value = "Ñuble"
$ws.Cells.Item(362, 4).Value = 45021
$ws.Cells.Item(362, 5).Value = 16
$ws.Cells.Item(362, 6).Value = 100112003
$ws.Cells.Item(362, 7).Value = "Ajo"
$ws.Cells.Item(362, 8).Value = "Chino"
$ws.Cells.Item(362, 9).Value = "1a (cosecha)"
$ws.Cells.Item(362, 10).Value = 50
$ws.Cells.Item(362, 11).Value = 18000
$ws.Cells.Item(362, 12).Value = 18000
$ws.Cells.Item(362, 13).Value = 18000
$ws.Cells.Item(362, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(362, 15).Value = "China"
$ws.Cells.Item(362, 16).Value = 1800
$ws.Cells.Item(362, 17).Value = 10
$ws.Cells.Item(362, 18).Value = "Hortaliza"
